$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column AB (28th column) from width 8 to width 7 (stored OOXML width).
# The ColumnWidth COM property is offset from the stored <col width> by ~0.83
# for this sheet's default font, so 7 - 0.83 = 6.17 yields a stored width of 7.
$ws.Columns.Item(28).ColumnWidth = 6.17

# Update row 5 values to their rounded (2 decimal place) figures.
$ws.Range("C5").Value = 16.52
$ws.Range("D5").Value = 0.6
$ws.Range("E5").Value = 47.74
$ws.Range("F5").Value = 39.37
$ws.Range("G5").Value = 17.27
$ws.Range("I5").Value = 26.81
$ws.Range("J5").Value = 12.17
$ws.Range("K5").Value = 18.24
$ws.Range("L5").Value = 19.51
$ws.Range("M5").Value = 20.33
$ws.Range("N5").Value = 5.56
$ws.Range("P5").Value = 24.67
$ws.Range("Q5").Value = 14.34
$ws.Range("R5").Value = 0.25
$ws.Range("S5").Value = 0.73
$ws.Range("T5").Value = 255.77
$ws.Range("U5").Value = 48.37
$ws.Range("V5").Value = 15.85
$ws.Range("W5").Value = 32.46
$ws.Range("X5").Value = 17.51
$ws.Range("Y5").Value = 2.27
$ws.Range("Z5").Value = 32.9
$ws.Range("AA5").Value = 14.08
$ws.Range("AB5").Value = 12.79
$ws.Range("AC5").Value = 15.04
$ws.Range("AD5").Value = 20.52
$ws.Range("AE5").Value = 0.21
$ws.Range("AF5").Value = 61.19
$ws.Range("AG5").Value = 9.289999999999999
# B5, H5, O5 and AH5 are unchanged in the diff.

# Remove row 6 entirely (data now ends at row 5).
$ws.Rows.Item(6).Delete()
